# Lithuania A Lyga - base update (21-04-2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) A handful of existing match rows got re-sorted: the two rows that share
#    the same Date swap all of their match data (columns B..AC) while the
#    running "id" in column A stays attached to the row position.
# ---------------------------------------------------------------------------
function Swap-Rows($r1, $r2) {
    $v1 = $ws.Range("B$r1`:AC$r1").Value()
    $v2 = $ws.Range("B$r2`:AC$r2").Value()
    $ws.Range("B$r1`:AC$r1").Value = $v2
    $ws.Range("B$r2`:AC$r2").Value = $v1
}

Swap-Rows 26 27
Swap-Rows 50 51
Swap-Rows 89 90
Swap-Rows 102 104
Swap-Rows 136 137

# ---------------------------------------------------------------------------
# 2) New fixtures added. The previous last row (142 / row 144) moves down to
#    row 148 (its id becomes 146), and four freshly-played matches are
#    inserted as the new rows 144-147.
# ---------------------------------------------------------------------------

# Preserve the old row 144 -> goes to row 148, with id renumbered to 146
$ws.Range("A148").Value = 146
$ws.Range("B148").Value = 7862048
$ws.Range("C148").Value = "Lithuania A Lyga"
$ws.Range("D148").Value = "Lithuania A Lyga"
$ws.Range("E148").Value = 45403.51736111111
$ws.Range("F148").Value = "FK Zalgiris Vilnius"
$ws.Range("G148").Value = "FK Dziugas Telsiai"
$ws.Range("K148").Value = 1.333
$ws.Range("L148").Value = 5
$ws.Range("M148").Value = 6
$ws.Range("N148").Value = 1.4
$ws.Range("O148").Value = 4.75
$ws.Range("P148").Value = 5.5
$ws.Range("Q148").Value = -1.25
$ws.Range("R148").Value = 1.925
$ws.Range("S148").Value = 1.875
$ws.Range("T148").Value = 2.5
$ws.Range("U148").Value = 1.8
$ws.Range("V148").Value = 2
$ws.Range("W148").Value = 0
$ws.Range("X148").Value = 0
$ws.Range("Y148").Value = 0
$ws.Range("Z148").Value = 0
$ws.Range("AA148").Value = 0

# New row 144 (id 142) - overwrite what used to be there
$ws.Range("A144").Value = 142
$ws.Range("B144").Value = 7862047
$ws.Range("C144").Value = "Lithuania A Lyga"
$ws.Range("D144").Value = "Lithuania A Lyga"
$ws.Range("E144").Value = 45402.375
$ws.Range("F144").Value = "FK Kauno Zalgiris"
$ws.Range("G144").Value = "Suduva Marijampole"
$ws.Range("H144").Value = 2
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = "H"
$ws.Range("K144").Value = 1.444
$ws.Range("L144").Value = 4
$ws.Range("M144").Value = 6
$ws.Range("N144").Value = 1.615
$ws.Range("O144").Value = 3.6
$ws.Range("P144").Value = 4.75
$ws.Range("Q144").Value = -0.75
$ws.Range("R144").Value = 1.825
$ws.Range("S144").Value = 1.975
$ws.Range("T144").Value = 2.25
$ws.Range("U144").Value = 2
$ws.Range("V144").Value = 1.8
$ws.Range("W144").Value = 0.615
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = 0.4125
$ws.Range("AA144").Value = -0.5
$ws.Range("AB144").Value = 1
$ws.Range("AC144").Value = -1

# New row 145 (id 143)
$ws.Range("A145").Value = 143
$ws.Range("B145").Value = 7862927
$ws.Range("C145").Value = "Lithuania A Lyga"
$ws.Range("D145").Value = "Lithuania A Lyga"
$ws.Range("E145").Value = 45402.5
$ws.Range("F145").Value = "FK Siauliai"
$ws.Range("G145").Value = "FK Dainava Alytus"
$ws.Range("H145").Value = 1
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = "H"
$ws.Range("K145").Value = 2
$ws.Range("L145").Value = 3.25
$ws.Range("M145").Value = 3.25
$ws.Range("N145").Value = 1.65
$ws.Range("O145").Value = 3.75
$ws.Range("P145").Value = 4.2
$ws.Range("Q145").Value = -0.75
$ws.Range("R145").Value = 1.9
$ws.Range("S145").Value = 1.9
$ws.Range("T145").Value = 2
$ws.Range("U145").Value = 1.775
$ws.Range("V145").Value = 2.025
$ws.Range("W145").Value = 0.6499999999999999
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 0.45
$ws.Range("AA145").Value = -0.5
$ws.Range("AB145").Value = -1
$ws.Range("AC145").Value = 1.025

# New row 146 (id 144)
$ws.Range("A146").Value = 144
$ws.Range("B146").Value = 7862928
$ws.Range("C146").Value = "Lithuania A Lyga"
$ws.Range("D146").Value = "Lithuania A Lyga"
$ws.Range("E146").Value = 45403.29166666666
$ws.Range("F146").Value = "Panevezys"
$ws.Range("G146").Value = "FK Transinvest"
$ws.Range("H146").Value = 1
$ws.Range("I146").Value = 1
$ws.Range("J146").Value = "D"
$ws.Range("K146").Value = 1.8
$ws.Range("L146").Value = 3.4
$ws.Range("M146").Value = 3.8
$ws.Range("N146").Value = 1.6
$ws.Range("O146").Value = 3.5
$ws.Range("P146").Value = 5
$ws.Range("Q146").Value = -0.75
$ws.Range("R146").Value = 1.825
$ws.Range("S146").Value = 1.975
$ws.Range("T146").Value = 2.25
$ws.Range("U146").Value = 1.95
$ws.Range("V146").Value = 1.85
$ws.Range("W146").Value = -1
$ws.Range("X146").Value = 2.5
$ws.Range("Y146").Value = -1
$ws.Range("Z146").Value = -1
$ws.Range("AA146").Value = 0.9750000000000001
$ws.Range("AB146").Value = -0.5
$ws.Range("AC146").Value = 0.425

# New row 147 (id 145)
$ws.Range("A147").Value = 145
$ws.Range("B147").Value = 7862929
$ws.Range("C147").Value = "Lithuania A Lyga"
$ws.Range("D147").Value = "Lithuania A Lyga"
$ws.Range("E147").Value = 45403.375
$ws.Range("F147").Value = "Hegelmann Litauen"
$ws.Range("G147").Value = "Banga Gargzdai"
$ws.Range("H147").Value = 2
$ws.Range("I147").Value = 2
$ws.Range("J147").Value = "D"
$ws.Range("K147").Value = 1.4
$ws.Range("L147").Value = 4.5
$ws.Range("M147").Value = 5.5
$ws.Range("N147").Value = 1.5
$ws.Range("O147").Value = 4.333
$ws.Range("P147").Value = 5
$ws.Range("Q147").Value = -1
$ws.Range("R147").Value = 1.825
$ws.Range("S147").Value = 1.975
$ws.Range("T147").Value = 2.75
$ws.Range("U147").Value = 2
$ws.Range("V147").Value = 1.8
$ws.Range("W147").Value = -1
$ws.Range("X147").Value = 3.333
$ws.Range("Y147").Value = -1
$ws.Range("Z147").Value = -1
$ws.Range("AA147").Value = 0.9750000000000001
$ws.Range("AB147").Value = 1
$ws.Range("AC147").Value = -1

# ---------------------------------------------------------------------------
# Match the "id"-column (A) and "Date"-column (E) formatting used by every
# other data row (bold/centered/bordered id, yyyy-mm-dd hh:mm:ss date) by
# cloning the format already used on row 2 (xlPasteFormats = -4122).
# ---------------------------------------------------------------------------
foreach ($r in 144..148) {
    $ws.Range("A2").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("E2").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)
}
